# Weekly price-data refresh: a new week's record is inserted at the top of
# this subset's data block (row 284), pushing the existing rows down by one.
# (Matches commit message "Fruta / hortaliza, semanal".)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 284; rows 284:306 shift down to 285:307.
$ws.Rows(284).Insert()

# Populate the new row 284 with this week's record.
$ws.Cells.Item(284, 1).Value  = 4
$ws.Cells.Item(284, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(284, 3).Value  = "Los Lagos"
$ws.Cells.Item(284, 4).Value  = (Get-Date -Year 2022 -Month 7 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(284, 5).Value  = 10
$ws.Cells.Item(284, 6).Value  = 100112037
$ws.Cells.Item(284, 7).Value  = "Cebollín"
$ws.Cells.Item(284, 8).Value  = "Sin especificar"
$ws.Cells.Item(284, 9).Value  = "Primera"
$ws.Cells.Item(284, 10).Value = 35
$ws.Cells.Item(284, 11).Value = 10000
$ws.Cells.Item(284, 12).Value = 10000
$ws.Cells.Item(284, 13).Value = 10000
$ws.Cells.Item(284, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(284, 15).Value = "Región Metropolitana"
$ws.Cells.Item(284, 16).Value = 278
$ws.Cells.Item(284, 17).Value = 36
$ws.Cells.Item(284, 18).Value = "Hortaliza"
